# LAPR5-146 #node #implementation Changes
#
# The peer-evaluation matrix (Grupo 051) on sheet "Fatura" scored every
# student pair with a raw numeric grade of 19. This replaces every one of
# those numeric self/peer-assessment scores with the qualitative label
# "Fair", leaving the diagonal (self-assessment) cells untouched (they
# stay blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gradeCells = @(
    "D4", "F4", "G4", "H4",
    "D5", "E5", "G5", "H5",
    "D6", "E6", "F6", "H6",
    "D7", "E7", "F7", "G7"
)

foreach ($addr in $gradeCells) {
    $ws.Range($addr).Value = "Fair"
}

# Move the active selection, matching the cursor position left behind by
# the author after making the edit (cosmetic, no data impact).
[void]$ws.Range("I10").Select()
